# Update column G ("K" = Strike#s -> K) values for rows 2-17 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 2
    9  = 2
    10 = 1
    11 = 1
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 2
    17 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
